$wb = $excel.ActiveWorkbook

# Rename sheets (task order identifiers updated)
$wb.Worksheets.Item(1).Name = "GNG_TO-1651168709293581"
$wb.Worksheets.Item(2).Name = "NB_TO-16511687126773167"
$wb.Worksheets.Item(3).Name = "RS_TO-1651168712679236"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511687127397335"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511687128156505"

# Sheet 1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16511687092515771.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687092765791.csv"
$ws1.Range("B4").Value = "go_stims-16511687092785795.csv"
$ws1.Range("B5").Value = "GNG_stims-1651168709292577.csv"

# Sheet 2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16511687116331182.csv"
$ws2.Range("B3").Value = "OB-16511687097215207.csv"
$ws2.Range("B4").Value = "ZB-match_0-16511687096469696.csv"
$ws2.Range("B5").Value = "OB-16511687109330406.csv"
$ws2.Range("B6").Value = "TB-16511687118523893.csv"
$ws2.Range("B7").Value = "OB-1651168710704904.csv"
$ws2.Range("B8").Value = "ZB-match_0-16511687095087662.csv"
$ws2.Range("B9").Value = "TB-16511687126514444.csv"
$ws2.Range("B10").Value = "ZB-match_9-16511687093705802.csv"

# Sheet 4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1651168712707984.csv"
$ws4.Range("B3").Value = "ZM_stims-165116871268223.csv"
$ws4.Range("B4").Value = "MM_stims-1651168712723266.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687127089481.csv"
$ws4.Range("B6").Value = "MM_stims-16511687127387373.csv"
$ws4.Range("B7").Value = "ZM_stims-1651168712724254.csv"

# Sheet 5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16511687127850873.csv"
$ws5.Range("B3").Value = "SAT_stims-16511687127694857.csv"
$ws5.Range("B4").Value = "vSAT_stims-16511687128000188.csv"
$ws5.Range("B5").Value = "SAT_stims-16511687127437363.csv"
